# Applies: row deletions in "settle_misc" and "settle_nil" sheets,
# selection/view changes on "settle_misc", "settle_nil" and "on_recon",
# and makes "settle_misc" the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- settle_misc (sheet1): delete the collector_app_id / receiver_app_id rows (old rows 15-16) ---
$wsMisc = $wb.Worksheets.Item("settle_misc")
$wsMisc.Rows.Item(15).Resize(2).Delete() | Out-Null

# --- settle_nil (sheet3): delete the same two attribute rows, plus everything
#     after the settlement.type row (old rows 18-31) ---
$wsNil = $wb.Worksheets.Item("settle_nil")
$wsNil.Rows.Item(15).Resize(2).Delete() | Out-Null
$wsNil.Rows.Item(16).Resize(14).Delete() | Out-Null

# --- View / selection updates ---

# settle_nil: scroll + new selection.
$wsNil.Activate()
$excel.ActiveWindow.ScrollRow = 8
$wsNil.Range("A17").Select() | Out-Null

# on_recon: scroll changes, selection stays on E18.
$wsRecon = $wb.Worksheets.Item("on_recon")
$wsRecon.Activate()
$excel.ActiveWindow.ScrollRow = 22
$wsRecon.Range("E18").Select() | Out-Null

# settle_misc becomes the active/selected sheet (saved as the active tab),
# with its scroll position and new selection.
$wsMisc.Activate()
$excel.ActiveWindow.ScrollRow = 12
$wsMisc.Range("A23").Select() | Out-Null
